$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''42.790.93'
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").Value = '''2.530.61'
$ws.Range("E3").Value = '  +0.53%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''315.62'
$ws.Range("E5").Value = '  +1.00%  '
$ws.Range("D6").Value = '''95.91'
$ws.Range("E6").Value = '  -0.12%  '
$ws.Range("D7").Value = '''0.576'
$ws.Range("E7").Value = '  -1.85%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = '''0.533'
$ws.Range("E9").Value = '  -1.13%  '
$ws.Range("D10").Value = '''36.15'
$ws.Range("E10").Value = '  -0.03%  '
$ws.Range("D11").Value = '''0.0810'
$ws.Range("E11").Value = '  -0.53%  '
$ws.Range("D12").Value = '''7.55'
$ws.Range("E12").Value = '  -0.90%  '
$ws.Range("E13").Value = '  -3.40%  '
$ws.Range("D14").Value = '''2.922.66'
$ws.Range("E14").Value = '  +0.96%  '
$ws.Range("D15").Value = '''2.562.37'
$ws.Range("E15").Value = '  +2.37%  '
$ws.Range("D16").Value = '''15.22'
$ws.Range("E16").Value = '  -2.23%  '
$ws.Range("D17").Value = '''0.854'
$ws.Range("E17").Value = '  -0.54%  '
$ws.Range("D18").Value = '''42.892.36'
$ws.Range("E18").Value = '  +0.80%  '
$ws.Range("D19").Value = '''6.80'
$ws.Range("E19").Value = '  +4.56%  '
$ws.Range("D20").Value = '''12.84'
$ws.Range("E20").Value = '  -0.92%  '
$ws.Range("D21").Value = '''0.0₃0964'
$ws.Range("E21").Value = '  -1.23%  '
$ws.Range("D22").Value = '''69.93'
$ws.Range("E22").Value = '  -2.31%  '
$ws.Range("D23").Value = '''253.04'
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("E24").Value = '  -0.26%  '
$ws.Range("D25").Value = '''2.07'
$ws.Range("E25").Value = '  +1.80%  '
$ws.Range("D26").Value = '''26.78'
$ws.Range("E26").Value = '  -1.15%  '
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("E28").Value = '  +3.27%  '
$ws.Range("D29").Value = '''40.60'
$ws.Range("E29").Value = '  +7.16%  '
$ws.Range("D30").Value = '''10.36'
$ws.Range("E30").Value = '  +2.11%  '
$ws.Range("E31").Value = '  -0.24%  '
$ws.Range("D32").Value = '''157.84'
$ws.Range("E32").Value = '  +1.60%  '
$ws.Range("D33").Value = '''2.17'
$ws.Range("E33").Value = '  +4.59%  '
$ws.Range("E34").Value = '  +0.14%  '
$ws.Range("D35").Value = '''19.12'
$ws.Range("E35").Value = '  -0.99%  '
$ws.Range("E36").Value = '  +2.01%  '
$ws.Range("D37").Value = '''0.0781'
$ws.Range("E37").Value = '  -0.72%  '
$ws.Range("E38").Value = '  -1.06%  '
$ws.Range("E39").Value = '  -1.27%  '
$ws.Range("D40").Value = '''23.51'
$ws.Range("E40").Value = '  -3.62%  '
$ws.Range("D41").Value = '''2.32'
$ws.Range("E41").Value = '  +14.64%  '
$ws.Range("D42").Value = '''3.83'
$ws.Range("E42").Value = '  -0.82%  '
$ws.Range("E43").Value = '  +0.36%  '
$ws.Range("E44").Value = '  +0.40%  '
$ws.Range("D45").Value = '''3.31'
$ws.Range("E45").Value = '  -2.06%  '
$ws.Range("D46").Value = '''2.039.15'
$ws.Range("E46").Value = '  +0.64%  '
$ws.Range("D47").Value = '''84.95'
$ws.Range("E47").Value = '  +0.48%  '
$ws.Range("D48").Value = '''9.02'
$ws.Range("E48").Value = '  +0.77%  '
$ws.Range("D49").Value = '''106.81'
$ws.Range("E49").Value = '  +5.17%  '
$ws.Range("D50").Value = '''74.97'
$ws.Range("E50").Value = '  +2.28%  '
$ws.Range("D51").Value = '''2.776.33'
$ws.Range("E51").Value = '  +0.81%  '
